# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-05 (serial 45204) to 2023-10-06 (serial 45205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 427 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45205
